# Junction_Flooding_457.xlsx : custom accuracy + 데이터 1000개
#
# - Updates simulation result values in rows 2-5 (new run of the model)
# - Removes the now-unused row 6 (dataset shrank from 5 data rows to 4)
# - Bumps several data-column widths from 7 -> 8 (and column T from 8 -> 9)
#   chars to better fit the new (wider) numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. New data values for rows 2-5 (columns A.. AH)
# ---------------------------------------------------------------------------

$headerCols = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH")

$row2 = @(45138.50694444445,15.854,10.494,3.735,33.762,26.097,12.381,37.784,19.197,7.826,11.635,13.329,13.955,3.981,12.407,17.155,10.858,3.133,2.01,181.475,34.412,11.452,22.298,11.285,3.17,19.349,10.115,9.106,10.918,14.036,3.313,33.905,6.108,14.317)
$row3 = @(45138.51388888889,12.491,8.699,1.645,27.032,21.31,9.747,37.98,15.125,6.45,9.393000000000001,10.785,11.338,3.141,9.775,13.7,8.571,1.445,0.954,141.466,27.379,9.023,17.975,9.304,1.99,18.511,7.97,7.212,8.532999999999999,11.335,1.246,34.909,4.887,11.281)
$row4 = @(45138.52083333334,3.843,2.436,0.887,8.302,6.126,2.952,15.843,4.654,1.917,2.578,3.31,3.469,0.976,3.008,4.181,2.848,0.907,0.457,38.501,8.666,2.776,5.53,2.778,0.855,7.256,2.452,2.325,2.743,3.469,0.766,14.914,1.397,3.473)
$row5 = @(45138.52777777778,1.92,1.11,0.61,4.15,2.86,1.45,8.34,2.33,0.91,1.14,1.66,1.74,0.49,1.5,2.05,1.53,0.67,0.3,15.59,4.36,1.39,2.69,1.34,0.54,3.71,1.23,1.22,1.43,1.73,0.55,7.85,0.65,1.74)

$dataRows = @{ 2 = $row2; 3 = $row3; 4 = $row4; 5 = $row5 }

foreach ($r in 2..5) {
    $values = $dataRows[$r]
    for ($i = 0; $i -lt $headerCols.Length; $i++) {
        $ws.Range("$($headerCols[$i])$r").Value = $values[$i]
    }
}

# ---------------------------------------------------------------------------
# 2. Drop row 6 - the refreshed dataset only has 4 data rows now
# ---------------------------------------------------------------------------

$ws.Rows.Item(6).Delete()

# ---------------------------------------------------------------------------
# 3. Widen a handful of data columns (7 -> 8 chars), and T (8 -> 9 chars)
# ---------------------------------------------------------------------------

$widerTo8 = @(2,3,6,7,11,12,13,15,17,22,24,27,29,30,34)
foreach ($c in $widerTo8) {
    $ws.Columns.Item($c).ColumnWidth = 7.166666666666668
}

$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
